$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A8").Value = "TestAutomation5oct"
$ws.Range("B8").Value = "TestAutomation5oct"
$ws.Range("C8").Value = "Facility_POC5oct"
$ws.Range("D8").Value = "Facility_POC5oct"
$ws.Range("E8").Value = "Pharmacy_POC5oct"
$ws.Range("F8").Value = "Pharmacy_POC5oct"
$ws.Range("H8").Value = "Epic101205"
